$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for rows 2..11 (columns D, L, M, N, O, P, Q, S, T)
$data = @(
    @{Row=2;  D=44441; L="Primera"; M=100; N=29000; O=30000; P=29500; Q="$/caja 12 kilos";    S=2458; T=12}
    @{Row=3;  D=44524; L="Primera"; M=200; N=23000; O=24000; P=23500; Q="$/caja 12 kilos";    S=1958; T=12}
    @{Row=4;  D=44475; L="Especial"; M=200; N=32000; O=33000; P=32500; Q="$/caja 12 kilos";    S=2708; T=12}
    @{Row=5;  D=44468; L="Primera"; M=200; N=29000; O=30000; P=29500; Q="$/bandeja 10 kilos"; S=2950; T=10}
    @{Row=6;  D=44496; L="Primera"; M=200; N=23000; O=24000; P=23500; Q="$/caja 12 kilos";    S=1958; T=12}
    @{Row=7;  D=44482; L="Primera"; M=160; N=25000; O=26000; P=25500; Q="$/caja 12 kilos";    S=2125; T=12}
    @{Row=8;  D=44167; L="Segunda"; M=200; N=18000; O=19000; P=18500; Q="$/caja 13 kilos";    S=1423; T=13}
    @{Row=9;  D=44545; L="Primera"; M=200; N=23000; O=24000; P=23500; Q="$/bandeja 12 kilos"; S=1958; T=12}
    @{Row=10; D=44489; L="Primera"; M=200; N=24000; O=25000; P=24500; Q="$/caja 12 kilos";    S=2042; T=12}
    @{Row=11; D=44160; L="Segunda"; M=200; N=19000; O=20000; P=19500; Q="$/caja 13 kilos";    S=1500; T=13}
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 4).Value = $rowData.D    # D: Fecha
    $ws.Cells.Item($r, 12).Value = $rowData.L   # L: Calidad
    $ws.Cells.Item($r, 13).Value = $rowData.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $rowData.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $rowData.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $rowData.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $rowData.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($r, 19).Value = $rowData.S   # S: Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $rowData.T   # T: Kg / unidad
}
